$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B29: it was stored as a text "4"; convert it to a real number 4.
$ws.Range("B29").Value = 4

# Add new row 30 with annotation data for Ying Tang.
$ws.Range("A30").Value = "Ying Tang"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "3"
$ws.Range("C30").Value = "elaborate"
$ws.Range("D30").Value = "ACK"
$ws.Range("E30").Value = "WRI"
$ws.Range("F30").Value = "9b81a0cf-ae6f-4476-b619-1b75e1becf94"
$ws.Range("G30").Value = "B1ae1lZRb_annotated.xlsx"
$ws.Range("H30").Value = "We will elaborate on this aspect in the final version of the paper."
